$wb = $excel.ActiveWorkbook

# --- Sheet 1: LH-TC-REGISTERATION-Reviews ---
$ws1 = $wb.Worksheets.Item("LH-TC-REGISTERATION-Reviews")

# Bring all review versions in line with the closed-out V1.0 registration review
$ws1.Range("E4").Value = "V1.0"
$ws1.Range("E5").Value = "V1.0"
$ws1.Range("E6").Value = "V1.0"

# Mark all reviewer verifications as Closed now that the review is done
$ws1.Range("J2").Value = "Closed"
$ws1.Range("J3").Value = "Closed"
$ws1.Range("J4").Value = "Closed"
$ws1.Range("J5").Value = "Closed"
$ws1.Range("J6").Value = "Closed"

$ws1.Range("J6").Select()

# --- Sheet 2: Version History ---
$ws2 = $wb.Worksheets.Item("Version History")

$ws2.Range("B2").Value = "Eman"
$ws2.Range("C2").Value = "intial review for the registeration feature"
$ws2.Range("D2").Value = 45767

$ws2.Range("A3").Value = "v1.1"
$ws2.Range("B3").Value = "Omar "
$ws2.Range("C3").Value = "update owner status for the reviews"
$ws2.Range("D3").Value = 45769

$ws2.Range("A4").Value = "V1.2"
$ws2.Range("B4").Value = "Eman"
$ws2.Range("C4").Value = "Verfiy thu updates`nclose the review"
$ws2.Range("D4").Value = 45769

$ws2.Range("B3").Select()
